# Fix the "Sim.  Parameter" (double space) typo to "Sim. Parameter" (single
# space) in the QOIList column (U) of Sheet1. This corrects the R list()
# literal strings used elsewhere to build the simulator's QOI dropdown.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldShort = 'list("Predicted Values", "Expected Values", "Sim.  Parameter")'
$newShort = 'list("Predicted Values", "Expected Values", "Sim. Parameter")'

$oldLong = 'list("Predicted Values", "Expected Values", "Probability Y > 1", "Sim.  Parameter")'
$newLong = 'list("Predicted Values", "Expected Values", "Probability Y > 1", "Sim. Parameter")'

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 21)  # column U
    $val = $cell.Value2
    if ($val -eq $oldShort) {
        $cell.Value = $newShort
    } elseif ($val -eq $oldLong) {
        $cell.Value = $newLong
    }
}
